$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 235.2675016666667
$ws.Range("H2").Value = 705.802505
$ws.Range("I2").Value = 0.5738994362335403
$ws.Range("J2").Value = 0.5738994362335402
$ws.Range("M2").Value = 159.4836373333333
$ws.Range("N2").Value = 478.450912
$ws.Range("O2").Value = 0.2983285084902258
$ws.Range("P2").Value = 0.2983285084902258
$ws.Range("Q2").Value = 37521.31691212607
$ws.Range("R2").Value = 337691.8522091346
$ws.Range("S2").Value = 0.1712105628349335
$ws.Range("T2").Value = 0.1712105628349335
$ws.Range("G3").Value = 235.2675016666667
$ws.Range("H3").Value = 705.802505
$ws.Range("I3").Value = 0.5738994362335403
$ws.Range("J3").Value = 0.5738994362335402
$ws.Range("O3").Value = 0.3227862111630279
$ws.Range("P3").Value = 0.3227862111630279
$ws.Range("Q3").Value = 40597.406480545
$ws.Range("R3").Value = 365376.6583249049
$ws.Range("S3").Value = 0.1852468246104222
$ws.Range("T3").Value = 0.1852468246104222
$ws.Range("G4").Value = 235.2675016666667
$ws.Range("H4").Value = 705.802505
$ws.Range("I4").Value = 0.5738994362335403
$ws.Range("J4").Value = 0.5738994362335402
$ws.Range("M4").Value = 74.38770566666666
$ws.Range("N4").Value = 223.163117
$ws.Range("O4").Value = 0.1391489036280481
$ws.Range("P4").Value = 0.1391489036280482
$ws.Range("Q4").Value = 17501.00966691201
$ws.Range("R4").Value = 157509.0870022081
$ws.Range("S4").Value = 0.07985747734465205
$ws.Range("T4").Value = 0.07985747734465205
$ws.Range("G5").Value = 235.2675016666667
$ws.Range("H5").Value = 705.802505
$ws.Range("I5").Value = 0.5738994362335403
$ws.Range("J5").Value = 0.5738994362335402
$ws.Range("M5").Value = 58.41461433333333
$ws.Range("N5").Value = 175.243843
$ws.Range("O5").Value = 0.1092697975759847
$ws.Range("P5").Value = 0.1092697975759848
$ws.Range("Q5").Value = 13743.06037502519
$ws.Range("R5").Value = 123687.5433752267
$ws.Range("S5").Value = 0.06270987522621072
$ws.Range("T5").Value = 0.06270987522621072
$ws.Range("G6").Value = 235.2675016666667
$ws.Range("H6").Value = 705.802505
$ws.Range("I6").Value = 0.5738994362335403
$ws.Range("J6").Value = 0.5738994362335402
$ws.Range("M6").Value = 69.746216
$ws.Range("N6").Value = 209.238648
$ws.Range("O6").Value = 0.1304665791427133
$ws.Range("P6").Value = 0.1304665791427133
$ws.Range("Q6").Value = 16409.01798902369
$ws.Range("R6").Value = 147681.1619012133
$ws.Range("S6").Value = 0.07487469621732175
$ws.Range("T6").Value = 0.07487469621732175
$ws.Range("I7").Value = 0.3286113026040369
$ws.Range("J7").Value = 0.3286113026040369
$ws.Range("M7").Value = 159.4836373333333
$ws.Range("N7").Value = 478.450912
$ws.Range("O7").Value = 0.2983285084902258
$ws.Range("P7").Value = 0.2983285084902258
$ws.Range("Q7").Value = 21484.47628182568
$ws.Range("R7").Value = 193360.2865364311
$ws.Range("S7").Value = 0.09803411977889259
$ws.Range("T7").Value = 0.09803411977889258
$ws.Range("I8").Value = 0.3286113026040369
$ws.Range("J8").Value = 0.3286113026040369
$ws.Range("O8").Value = 0.3227862111630279
$ws.Range("P8").Value = 0.3227862111630279
$ws.Range("S8").Value = 0.1060711973129043
$ws.Range("T8").Value = 0.1060711973129043
$ws.Range("I9").Value = 0.3286113026040369
$ws.Range("J9").Value = 0.3286113026040369
$ws.Range("M9").Value = 74.38770566666666
$ws.Range("N9").Value = 223.163117
$ws.Range("O9").Value = 0.1391489036280481
$ws.Range("P9").Value = 0.1391489036280482
$ws.Range("Q9").Value = 10020.97095838494
$ws.Range("R9").Value = 90188.73862546447
$ws.Range("S9").Value = 0.0457259024771365
$ws.Range("T9").Value = 0.0457259024771365
$ws.Range("I10").Value = 0.3286113026040369
$ws.Range("J10").Value = 0.3286113026040369
$ws.Range("M10").Value = 58.41461433333333
$ws.Range("N10").Value = 175.243843
$ws.Range("O10").Value = 0.1092697975759847
$ws.Range("P10").Value = 0.1092697975759848
$ws.Range("Q10").Value = 7869.192207683539
$ws.Range("R10").Value = 70822.72986915185
$ws.Range("S10").Value = 0.03590729051672378
$ws.Range("T10").Value = 0.03590729051672378
$ws.Range("I11").Value = 0.3286113026040369
$ws.Range("J11").Value = 0.3286113026040369
$ws.Range("M11").Value = 69.746216
$ws.Range("N11").Value = 209.238648
$ws.Range("O11").Value = 0.1304665791427133
$ws.Range("P11").Value = 0.1304665791427133
$ws.Range("Q11").Value = 9395.703211027158
$ws.Range("R11").Value = 84561.32889924443
$ws.Range("S11").Value = 0.0428727925183797
$ws.Range("T11").Value = 0.0428727925183797
$ws.Range("G12").Value = 0.325805
$ws.Range("H12").Value = 0.977415
$ws.Range("I12").Value = 0.0007947519504286909
$ws.Range("J12").Value = 0.0007947519504286907
$ws.Range("M12").Value = 159.4836373333333
$ws.Range("N12").Value = 478.450912
$ws.Range("O12").Value = 0.2983285084902258
$ws.Range("P12").Value = 0.2983285084902258
$ws.Range("Q12").Value = 51.96056646138668
$ws.Range("R12").Value = 467.64509815248
$ws.Range("S12").Value = 0.0002370971639910892
$ws.Range("T12").Value = 0.0002370971639910892
$ws.Range("G13").Value = 0.325805
$ws.Range("H13").Value = 0.977415
$ws.Range("I13").Value = 0.0007947519504286909
$ws.Range("J13").Value = 0.0007947519504286907
$ws.Range("O13").Value = 0.3227862111630279
$ws.Range("P13").Value = 0.3227862111630279
$ws.Range("Q13").Value = 56.220421115085
$ws.Range("R13").Value = 505.983790035765
$ws.Range("S13").Value = 0.0002565349708933037
$ws.Range("T13").Value = 0.0002565349708933037
$ws.Range("G14").Value = 0.325805
$ws.Range("H14").Value = 0.977415
$ws.Range("I14").Value = 0.0007947519504286909
$ws.Range("J14").Value = 0.0007947519504286907
$ws.Range("M14").Value = 74.38770566666666
$ws.Range("N14").Value = 223.163117
$ws.Range("O14").Value = 0.1391489036280481
$ws.Range("P14").Value = 0.1391489036280482
$ws.Range("Q14").Value = 24.23588644472833
$ws.Range("R14").Value = 218.122978002555
$ws.Range("S14").Value = 0.0001105888625584052
$ws.Range("T14").Value = 0.0001105888625584052
$ws.Range("G15").Value = 0.325805
$ws.Range("H15").Value = 0.977415
$ws.Range("I15").Value = 0.0007947519504286909
$ws.Range("J15").Value = 0.0007947519504286907
$ws.Range("M15").Value = 58.41461433333333
$ws.Range("N15").Value = 175.243843
$ws.Range("O15").Value = 0.1092697975759847
$ws.Range("P15").Value = 0.1092697975759848
$ws.Range("Q15").Value = 19.03177342287167
$ws.Range("R15").Value = 171.285960805845
$ws.Range("S15").Value = 0.00008684238474646213
$ws.Range("T15").Value = 0.00008684238474646211
$ws.Range("G16").Value = 0.325805
$ws.Range("H16").Value = 0.977415
$ws.Range("I16").Value = 0.0007947519504286909
$ws.Range("J16").Value = 0.0007947519504286907
$ws.Range("M16").Value = 69.746216
$ws.Range("N16").Value = 209.238648
$ws.Range("O16").Value = 0.1304665791427133
$ws.Range("P16").Value = 0.1304665791427133
$ws.Range("Q16").Value = 22.72366590388
$ws.Range("R16").Value = 204.51299313492
$ws.Range("S16").Value = 0.0001036885682394306
$ws.Range("T16").Value = 0.0001036885682394306
$ws.Range("G17").Value = 39.46134166666666
$ws.Range("H17").Value = 118.384025
$ws.Range("I17").Value = 0.09625996610278018
$ws.Range("J17").Value = 0.09625996610278018
$ws.Range("M17").Value = 159.4836373333333
$ws.Range("N17").Value = 478.450912
$ws.Range("O17").Value = 0.2983285084902258
$ws.Range("P17").Value = 0.2983285084902258
$ws.Range("Q17").Value = 6293.438303053422
$ws.Range("R17").Value = 56640.9447274808
$ws.Range("S17").Value = 0.0287170921147621
$ws.Range("T17").Value = 0.0287170921147621
$ws.Range("G18").Value = 39.46134166666666
$ws.Range("H18").Value = 118.384025
$ws.Range("I18").Value = 0.09625996610278018
$ws.Range("J18").Value = 0.09625996610278018
$ws.Range("O18").Value = 0.3227862111630279
$ws.Range("P18").Value = 0.3227862111630279
$ws.Range("Q18").Value = 6809.389807603474
$ws.Range("R18").Value = 61284.50826843127
$ws.Range("S18").Value = 0.03107138974499791
$ws.Range("T18").Value = 0.03107138974499791
$ws.Range("G19").Value = 39.46134166666666
$ws.Range("H19").Value = 118.384025
$ws.Range("I19").Value = 0.09625996610278018
$ws.Range("J19").Value = 0.09625996610278018
$ws.Range("M19").Value = 74.38770566666666
$ws.Range("N19").Value = 223.163117
$ws.Range("O19").Value = 0.1391489036280481
$ws.Range("P19").Value = 0.1391489036280482
$ws.Range("Q19").Value = 2935.438669111769
$ws.Range("R19").Value = 26418.94802200592
$ws.Range("S19").Value = 0.01339446874647494
$ws.Range("T19").Value = 0.01339446874647494
$ws.Range("G20").Value = 39.46134166666666
$ws.Range("H20").Value = 118.384025
$ws.Range("I20").Value = 0.09625996610278018
$ws.Range("J20").Value = 0.09625996610278018
$ws.Range("M20").Value = 58.41461433333333
$ws.Range("N20").Value = 175.243843
$ws.Range("O20").Value = 0.1092697975759847
$ws.Range("P20").Value = 0.1092697975759848
$ws.Range("Q20").Value = 2305.11905453423
$ws.Range("R20").Value = 20746.07149080807
$ws.Range("S20").Value = 0.01051830701072194
$ws.Range("T20").Value = 0.01051830701072195
$ws.Range("G21").Value = 39.46134166666666
$ws.Range("H21").Value = 118.384025
$ws.Range("I21").Value = 0.09625996610278018
$ws.Range("J21").Value = 0.09625996610278018
$ws.Range("M21").Value = 69.746216
$ws.Range("N21").Value = 209.238648
$ws.Range("O21").Value = 0.1304665791427133
$ws.Range("P21").Value = 0.1304665791427133
$ws.Range("Q21").Value = 2752.279259533133
$ws.Range("R21").Value = 24770.5133357982
$ws.Range("S21").Value = 0.01255870848582327
$ws.Range("T21").Value = 0.01255870848582327
$ws.Range("G22").Value = 0.178139
$ws.Range("H22").Value = 0.534417
$ws.Range("I22").Value = 0.0004345431092138444
$ws.Range("J22").Value = 0.0004345431092138443
$ws.Range("M22").Value = 159.4836373333333
$ws.Range("N22").Value = 478.450912
$ws.Range("O22").Value = 0.2983285084902258
$ws.Range("P22").Value = 0.2983285084902258
$ws.Range("Q22").Value = 28.41025567092267
$ws.Range("R22").Value = 255.692301038304
$ws.Range("S22").Value = 0.0001296365976464715
$ws.Range("T22").Value = 0.0001296365976464715
$ws.Range("G23").Value = 0.178139
$ws.Range("H23").Value = 0.534417
$ws.Range("I23").Value = 0.0004345431092138444
$ws.Range("J23").Value = 0.0004345431092138443
$ws.Range("O23").Value = 0.3227862111630279
$ws.Range("P23").Value = 0.3227862111630279
$ws.Range("Q23").Value = 30.739398097083
$ws.Range("R23").Value = 276.654582873747
$ws.Range("S23").Value = 0.0001402645238101387
$ws.Range("T23").Value = 0.0001402645238101387
$ws.Range("G24").Value = 0.178139
$ws.Range("H24").Value = 0.534417
$ws.Range("I24").Value = 0.0004345431092138444
$ws.Range("J24").Value = 0.0004345431092138443
$ws.Range("M24").Value = 74.38770566666666
$ws.Range("N24").Value = 223.163117
$ws.Range("O24").Value = 0.1391489036280481
$ws.Range("P24").Value = 0.1391489036280482
$ws.Range("Q24").Value = 13.25135149975433
$ws.Range("R24").Value = 119.262163497789
$ws.Range("S24").Value = 0.00006046619722622963
$ws.Range("T24").Value = 0.00006046619722622963
$ws.Range("G25").Value = 0.178139
$ws.Range("H25").Value = 0.534417
$ws.Range("I25").Value = 0.0004345431092138444
$ws.Range("J25").Value = 0.0004345431092138443
$ws.Range("M25").Value = 58.41461433333333
$ws.Range("N25").Value = 175.243843
$ws.Range("O25").Value = 0.1092697975759847
$ws.Range("P25").Value = 0.1092697975759848
$ws.Range("Q25").Value = 10.40592098272567
$ws.Range("R25").Value = 93.653288844531
$ws.Range("S25").Value = 0.00004748243758183581
$ws.Range("T25").Value = 0.00004748243758183581
$ws.Range("G26").Value = 0.178139
$ws.Range("H26").Value = 0.534417
$ws.Range("I26").Value = 0.0004345431092138444
$ws.Range("J26").Value = 0.0004345431092138443
$ws.Range("M26").Value = 69.746216
$ws.Range("N26").Value = 209.238648
$ws.Range("O26").Value = 0.1304665791427133
$ws.Range("P26").Value = 0.1304665791427133
$ws.Range("Q26").Value = 12.424521172024
$ws.Range("R26").Value = 111.820690548216
$ws.Range("S26").Value = 0.0001036885682394306
$ws.Range("T26").Value = 0.0001036885682394306
